$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells keep their original text representation
# (values such as "1.003" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "22.455.80"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.573.71"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "291.09"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.3740"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("D8").Value = "49.90"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.3401"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "0.07555"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "1.139"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "21.39"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "5.990"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "6.940"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "1.566.14"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "0.00001122"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "91.02"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "0.06750"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "6.265"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "16.40"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "22.446.21"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "2.337"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("D26").Value = "2.590"
$ws.Range("E26").Value = "  -5.81%  "
$ws.Range("D27").Value = "20.15"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "148.77"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").Value = "5.009"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "1.742.50"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").Value = "1.053"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").Value = "6.126"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").Value = "1.981"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "9.823"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").Value = "0.08415"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("D37").Value = "1.382"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("D38").Value = "0.02463"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").Value = "0.2291"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").Value = "0.06536"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").Value = "5.455"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "11.31"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").Value = "0.6261"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "13.99"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "3.810"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").Value = "0.5829"
$ws.Range("D48").Value = "2.087"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "129.19"
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").Value = "1.227"
$ws.Range("E50").Value = "  -5.33%  "
$ws.Range("D51").Value = "0.07324"
$ws.Range("E51").Value = "  -0.02%  "
